$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.63
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("X2").Value = 7.5
$ws.Range("AF2").Value = 81
